$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4103.75
$ws.Range("J64").Value = 3099.5
$ws.Range("L64").Value = 3099.5
$ws.Range("N64").Value = -3595.5

$ws.Range("H67").Value = 4103.75
$ws.Range("J67").Value = 3099.5
$ws.Range("L67").Value = 3099.5
$ws.Range("N67").Value = -4815.5

$ws.Range("H106").Value = 225433.69
$ws.Range("I106").Value = 294849.75
$ws.Range("J106").Value = 5616.1665
$ws.Range("K106").Value = 294849.75
$ws.Range("L106").Value = 5616.1665
$ws.Range("M106").Value = -294218.75
$ws.Range("N106").Value = -6878.1665

$ws.Range("H108").Value = 53120.625
$ws.Range("J108").Value = 53120.625
$ws.Range("L108").Value = 53120.625
$ws.Range("N108").Value = -60800.625

$ws.Range("H109").Value = 99922.5
$ws.Range("J109").Value = 99922.5
$ws.Range("L109").Value = 99922.5
$ws.Range("N109").Value = -102696.5

$ws.Range("H110").Value = 64640
$ws.Range("J110").Value = 64640
$ws.Range("L110").Value = 64640
$ws.Range("N110").Value = -72820

$ws.Range("H133").Value = 92754.13
$ws.Range("J133").Value = 92754.13
$ws.Range("L133").Value = 92754.13
$ws.Range("N133").Value = -102874.13

$ws.Range("H134").Value = 59718.125
$ws.Range("J134").Value = 68990.836
$ws.Range("L134").Value = 68990.836
$ws.Range("N134").Value = -79130.836

$ws.Range("H136").Value = 68273.86
$ws.Range("J136").Value = 81183.39999999999
$ws.Range("L136").Value = 81183.39999999999
$ws.Range("N136").Value = -91383.39999999999

$ws.Range("H137").Value = 447735.53
$ws.Range("I137").Value = 1264.7084
$ws.Range("J137").Value = 1787148
$ws.Range("K137").Value = 3794.1252
$ws.Range("L137").Value = 5361444
$ws.Range("M137").Value = -1244.1252
$ws.Range("N137").Value = -5366544

$ws.Range("H138").Value = 17650.092
$ws.Range("J138").Value = 2618.432
$ws.Range("L138").Value = 7855.295999999999
$ws.Range("N138").Value = -18135.296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11634.63
$ws.Range("I32").Value = 7304.923
$ws.Range("K32").Value = 7304.923
$ws.Range("M32").Value = -7017.923

$ws.Range("H52").Value = 54544
$ws.Range("J52").Value = 54544
$ws.Range("L52").Value = 54544
$ws.Range("N52").Value = -55180

$ws.Range("H104").Value = 44996
$ws.Range("J104").Value = 43745
$ws.Range("L104").Value = 43745
$ws.Range("N104").Value = -50733

$ws.Range("H122").Value = 3507.465
$ws.Range("I122").Value = 3669.8215
$ws.Range("J122").Value = 3204.4
$ws.Range("K122").Value = 11009.4645
$ws.Range("L122").Value = 9613.200000000001
$ws.Range("M122").Value = -8559.4645
$ws.Range("N122").Value = -14513.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 64930
$ws.Range("J6").Value = 64930
$ws.Range("L6").Value = 64930
$ws.Range("N6").Value = -65156

$ws.Range("H26").Value = 6667
$ws.Range("I26").Value = 6667
$ws.Range("K26").Value = 6667
$ws.Range("M26").Value = -6375

$ws.Range("H51").Value = 77648.664
$ws.Range("J51").Value = 77648.664
$ws.Range("L51").Value = 77648.664
$ws.Range("N51").Value = -78630.664

$ws.Range("H94").Value = 938.7826
$ws.Range("I94").Value = 992.1875
$ws.Range("J94").Value = 816.7143
$ws.Range("K94").Value = 992.1875
$ws.Range("L94").Value = 816.7143
$ws.Range("M94").Value = -541.1875
$ws.Range("N94").Value = -1718.7143

$ws.Range("H115").Value = 90467.75
$ws.Range("J115").Value = 99956.664
$ws.Range("L115").Value = 99956.664
$ws.Range("N115").Value = -103090.664

$ws.Range("H117").Value = 80392.42999999999
$ws.Range("J117").Value = 80392.42999999999
$ws.Range("L117").Value = 80392.42999999999
$ws.Range("N117").Value = -89570.42999999999

$ws.Range("H119").Value = 60992.25
$ws.Range("J119").Value = 60992.25
$ws.Range("L119").Value = 60992.25
$ws.Range("N119").Value = -70668.25

$ws.Range("H127").Value = 69883
$ws.Range("J127").Value = 69883
$ws.Range("L127").Value = 69883
$ws.Range("N127").Value = -79803

$ws.Range("H132").Value = 94996
$ws.Range("J132").Value = 94996
$ws.Range("L132").Value = 94996
$ws.Range("N132").Value = -105116

$ws.Range("H134").Value = 2088.2
$ws.Range("I134").Value = 1595.2222
$ws.Range("K134").Value = 4785.6666
$ws.Range("M134").Value = -2250.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1451.6666
$ws.Range("I105").Value = 1152.3572
$ws.Range("J105").Value = 2499.25
$ws.Range("K105").Value = 1152.3572
$ws.Range("L105").Value = 2499.25
$ws.Range("M105").Value = 594.6428000000001
$ws.Range("N105").Value = -5993.25

$ws.Range("H118").Value = 73102.39999999999
$ws.Range("J118").Value = 73102.39999999999
$ws.Range("L118").Value = 73102.39999999999
$ws.Range("N118").Value = -76416.39999999999

$ws.Range("H119").Value = 95192
$ws.Range("J119").Value = 95192
$ws.Range("L119").Value = 95192
$ws.Range("N119").Value = -104868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 640.55554
$ws.Range("I34").Value = 130
$ws.Range("J34").Value = 1442.8572
$ws.Range("K34").Value = 390
$ws.Range("L34").Value = 4328.571599999999
$ws.Range("M34").Value = -306
$ws.Range("N34").Value = -4496.571599999999

$ws.Range("H131").Value = 36198.863
$ws.Range("I131").Value = 77772.53999999999
$ws.Range("J131").Value = 2420.25
$ws.Range("K131").Value = 233317.62
$ws.Range("L131").Value = 7260.75
$ws.Range("M131").Value = -228277.62
$ws.Range("N131").Value = -17340.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 40658.465
$ws.Range("J109").Value = 40658.465
$ws.Range("L109").Value = 40658.465
$ws.Range("N109").Value = -42738.465

$ws.Range("H132").Value = 1795.9546
$ws.Range("I132").Value = 1278.3889
$ws.Range("K132").Value = 3835.1667
$ws.Range("M132").Value = -1305.1667

$ws.Range("H135").Value = 58329.832
$ws.Range("J135").Value = 58329.832
$ws.Range("L135").Value = 58329.832
$ws.Range("N135").Value = -68469.83199999999

$ws.Range("H140").Value = 68912.5
$ws.Range("J140").Value = 75550
$ws.Range("L140").Value = 75550
$ws.Range("N140").Value = -85910

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16383.4375
$ws.Range("I7").Value = 18048.691
$ws.Range("K7").Value = 18048.691
$ws.Range("M7").Value = -17936.691

$ws.Range("H42").Value = 4227.5
$ws.Range("J42").Value = 4227.5
$ws.Range("L42").Value = 4227.5
$ws.Range("N42").Value = -5353.5

$ws.Range("H49").Value = 4227.5
$ws.Range("J49").Value = 4227.5
$ws.Range("L49").Value = 4227.5
$ws.Range("N49").Value = -4521.5

$ws.Range("H68").Value = 3943.4614
$ws.Range("I68").Value = 4285.222
$ws.Range("J68").Value = 3174.5
$ws.Range("K68").Value = 4285.222
$ws.Range("L68").Value = 3174.5
$ws.Range("M68").Value = -3536.222
$ws.Range("N68").Value = -4672.5

$ws.Range("H71").Value = 3943.4614
$ws.Range("I71").Value = 4285.222
$ws.Range("J71").Value = 3174.5
$ws.Range("K71").Value = 21426.11
$ws.Range("L71").Value = 15872.5
$ws.Range("M71").Value = -17682.11
$ws.Range("N71").Value = -23360.5

$ws.Range("H118").Value = 95678
$ws.Range("J118").Value = 95678
$ws.Range("L118").Value = 95678
$ws.Range("N118").Value = -98992

$ws.Range("H126").Value = 16383.4375
$ws.Range("I126").Value = 18048.691
$ws.Range("K126").Value = 54146.073
$ws.Range("M126").Value = -51676.073

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 64008.43
$ws.Range("J121").Value = 64008.43
$ws.Range("L121").Value = 64008.43
$ws.Range("N121").Value = -67502.42999999999

$ws.Range("H126").Value = 13799.823
$ws.Range("I126").Value = 1887
$ws.Range("J126").Value = 24389
$ws.Range("K126").Value = 5661
$ws.Range("L126").Value = 73167
$ws.Range("M126").Value = -3191
$ws.Range("N126").Value = -78107
